function HexToLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b*65536 + $g*256 + $r
}

$p = $ppt.ActivePresentation

# The presentation's single reachable DrawingML theme (ppt/theme/theme2.xml,
# the one used by the slide master) currently holds the "Integral" palette.
# Swap its 12 theme colours for the stock "Office" palette (the palette that
# currently lives, unused, in ppt/theme/theme1.xml).
$tcs = $p.Slides.Item(1).ThemeColorScheme

$officeColors = @(
    "000000", # 1 dk1
    "FFFFFF", # 2 lt1
    "44546A", # 3 dk2
    "E7E6E6", # 4 lt2
    "5B9BD5", # 5 accent1
    "ED7D31", # 6 accent2
    "A5A5A5", # 7 accent3
    "FFC000", # 8 accent4
    "4472C4", # 9 accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $tcs.Colors($i).RGB = HexToLong $officeColors[$i - 1]
}

Write-Host "theme colours updated"
